$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates — same edits apply to both the "展览" sheet
# and the "全部类型" sheet (the latter aggregates all event types).
$updates = @{
    "F4"  = 12460
    "F5"  = 1276
    "F6"  = 145
    "F8"  = 88
    "F10" = 196
    "F16" = 376
    "F17" = 3913
    "F18" = 96
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
